# CV updated with solidity skills
#
# Slide 2, shape "TextBox 42" (work-experience bullet list):
#   - a new bullet "Researched and analyzed the EVM technology stack and early
#     Solidity codebase to write simple contracts." is inserted right before the
#     existing "Extensive research into Steem ..." bullet, styled like the
#     "Recognized as a leading ..." bullet (374151 / Söhne / b=0 / i=0).
#   - the "Extensive research into Steem ..." bullet gets an extra sentence
#     appended: " Wrote several contracts using C++; also written test
#     scripts; deployment using CLI."
#   - the textbox grows taller to fit the extra bullet, so its position/size
#     (and the position of the textbox right below it, "TextBox 47") shift.

# --- helpers: PowerPoint COM shape geometry is expressed in points, but the
# underlying OOXML stores EMU (1 pt = 12700 EMU). The host keeps shape
# position/size internally as 32-bit floats, so plugging target_emu/12700
# straight into .Top/.Height can truncate to target_emu-1. Solve for a point
# value that truncates back to exactly the desired EMU so the saved XML
# matches bit-for-bit. ---
function Get-EmuFromPt($pt) {
    $f32 = [float]$pt
    $f64 = [double]$f32
    $emu = $f64 * 12700.0
    return [math]::Floor($emu)
}

function Find-PtForEmu($targetEmu) {
    $basePt = $targetEmu / 12700.0
    for ($i = 0; $i -le 10000; $i++) {
        $candidate = $basePt + ($i * 0.0000001)
        if ((Get-EmuFromPt $candidate) -eq $targetEmu) {
            return $candidate
        }
    }
    for ($i = 1; $i -le 10000; $i++) {
        $candidate = $basePt - ($i * 0.0000001)
        if ((Get-EmuFromPt $candidate) -eq $targetEmu) {
            return $candidate
        }
    }
    throw "Find-PtForEmu: no point value found for target EMU $targetEmu"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shpExtras = $s.Shapes.Item(2)   # "TextBox 42"
$shpBelow  = $s.Shapes.Item(3)   # "TextBox 47"

$tr = $shpExtras.TextFrame.TextRange

# Paragraph 3 is currently "Extensive research into Steem blockchain's
# Tokenomics model & deep dived into EOS blockchain codebase."
$oldResearchPara = $tr.Paragraphs(3, 1)

# --- Insert the new "Researched and analyzed ..." bullet right before it ---
$newText = "Researched and analyzed the EVM technology stack and early Solidity codebase to write simple contracts."
$oldResearchPara.InsertBefore($newText + "`r")

# The newly inserted text is now paragraph 3; restyle it to match the
# "Recognized as a leading ..." bullet's look (374151 / Söhne / not bold / not italic).
$newPara = $tr.Paragraphs(3, 1)

$run1 = $newPara.Characters(1, 15)                                    # "Researched and "
$run2 = $newPara.Characters(16, 8)                                    # "analyzed"
$run3 = $newPara.Characters(24, $newPara.Length - 23)                 # " the EVM technology stack and early Solidity codebase to write simple contracts."

foreach ($run in @($run1, $run2, $run3)) {
    $run.Font.Name = "Söhne"
    $run.Font.Bold = $false
    $run.Font.Italic = $false
    $run.Font.Color.RGB = 0x514137   # OLE (BGR) form of srgbClr 374151
    $run.Font.Shadow = $false        # forces an <a:effectLst/> to be emitted, matching the target markup
}

# --- Append the extra sentence to the (now 4th) "Extensive research into Steem ..." bullet ---
$researchPara = $tr.Paragraphs(4, 1)
$tailRun = $researchPara.Characters(30, $researchPara.Length - 29)    # " blockchain's Tokenomics model & deep dived into EOS blockchain codebase."
$tailRun.Text = " blockchain" + [char]0x2019 + "s Tokenomics model & deep dived into EOS blockchain codebase. Wrote several contracts using C++; also written test scripts; deployment using CLI."

# --- Resize/reposition the textbox to fit the new bullet, and shift the box below it ---
$shpExtras.Top    = Find-PtForEmu 4780095
$shpExtras.Height = Find-PtForEmu 1546257

$shpBelow.Top = Find-PtForEmu 3747787
